# Prepend two new floating-shape (v:group) picture runs to the start of
# the document's single paragraph, per the target diff. The new content
# must land inside the SAME paragraph as the existing picture runs (no
# new paragraph mark introduced) and must precede all existing runs.
#
# This document's paragraph contains no visible text -- every run is a
# <w:pict> floating shape -- so Word's character-range addressing only
# exposes a single position (the paragraph mark). Collapsing the range
# to that trailing position and calling InsertXML with a full <w:p>
# fragment merges the new runs into the existing paragraph, ahead of
# the pre-existing (non-addressable) picture runs.

$d = $word.ActiveDocument

$newRunsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:pict><v:group id="_x0000_s1088" style="position:absolute;margin-left:274.5pt;margin-top:327.15pt;width:138pt;height:155.6pt;z-index:251751424" coordorigin="5230,2988" coordsize="2760,3112"><v:roundrect id="_x0000_s1089" style="position:absolute;left:5230;top:2988;width:2760;height:3112" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1089" inset=",1mm"><w:txbxContent><w:p><w:pPr><w:spacing w:before="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Mô hình nhận </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>dạng</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect><v:roundrect id="_x0000_s1090" style="position:absolute;left:5570;top:3691;width:2092;height:668;v-text-anchor:middle" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1090" inset=",0,,2mm"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>PCA</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect><v:roundrect id="_x0000_s1091" style="position:absolute;left:5570;top:4477;width:2092;height:668;v-text-anchor:middle" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1091" inset=",0,,2mm"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>LDA</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect><v:roundrect id="_x0000_s1092" style="position:absolute;left:5570;top:5277;width:2092;height:668;v-text-anchor:middle" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1092" inset=",0,,2mm"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>SFS</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect></v:group></w:pict></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:pict><v:group id="_x0000_s1083" style="position:absolute;margin-left:69.6pt;margin-top:327.15pt;width:138pt;height:155.6pt;z-index:251750400" coordorigin="5230,2988" coordsize="2760,3112"><v:roundrect id="_x0000_s1084" style="position:absolute;left:5230;top:2988;width:2760;height:3112" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1084" inset=",1mm"><w:txbxContent><w:p><w:pPr><w:spacing w:before="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Trích chọn đặc tính</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect><v:roundrect id="_x0000_s1085" style="position:absolute;left:5570;top:3691;width:2092;height:668;v-text-anchor:middle" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1085" inset=",0,,2mm"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>PCA</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect><v:roundrect id="_x0000_s1086" style="position:absolute;left:5570;top:4477;width:2092;height:668;v-text-anchor:middle" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1086" inset=",0,,2mm"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>LDA</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect><v:roundrect id="_x0000_s1087" style="position:absolute;left:5570;top:5277;width:2092;height:668;v-text-anchor:middle" arcsize="10923f"><v:textbox style="mso-next-textbox:#_x0000_s1087" inset=",0,,2mm"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>SFS</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect></v:group></w:pict></w:r></w:p>
'@

$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$insertionPoint.InsertXML($newRunsXml)
